# Auto-generated Excel COM-interop script
# Commit: Updated symbol list on Fri Jan 27 00:59:45 UTC 2023 with GitHub Actions
# Refreshes the crypto price table (columns B-G, rows 2-51) with the latest
# coinranking.com snapshot: prices, 1h volume %, rank reshuffles, and the
# "Data"/"Hora" timestamp columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as TEXT
# (the source feed always writes these columns as text/inline strings, e.g.
# "303.70", "-1.59 percent", "27-1-2023", "0" must stay text and not become
# numbers, percentages or dates). NumberFormat is restored to General
# afterwards so no stray text-style is left on the cell.
function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $range = $Worksheet.Range($CellRef)
    $range.NumberFormat = "@"
    $range.Value = $Text
    $range.NumberFormat = "General"
}

# Cell reference -> new text value, in document order (row 2 through row 51)
$updates = [ordered]@{
    "D2" = "303.70"
    "E2" = "-1.59%"
    "F2" = "27-1-2023"
    "G2" = "0"
    "D3" = "36.00"
    "E3" = "-1.00%"
    "F3" = "27-1-2023"
    "G3" = "0"
    "D4" = "5.020"
    "E4" = "-2.23%"
    "F4" = "27-1-2023"
    "G4" = "0"
    "D5" = "0.07993"
    "E5" = "-2.10%"
    "F5" = "27-1-2023"
    "G5" = "0"
    "D6" = "1.856"
    "E6" = "-5.41%"
    "F6" = "27-1-2023"
    "G6" = "0"
    "B7" = "GateToken"
    "C7" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D7" = "4.163"
    "E7" = "0.65%"
    "F7" = "27-1-2023"
    "G7" = "0"
    "B8" = "KuCoinToken"
    "C8" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "D8" = "7.831"
    "E8" = "0.93%"
    "F8" = "27-1-2023"
    "G8" = "0"
    "B9" = "MXToken"
    "C9" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D9" = "0.9291"
    "E9" = "-0.96%"
    "F9" = "27-1-2023"
    "G9" = "0"
    "B10" = "LiechtensteinCryptoassetsExchange"
    "C10" = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D10" = "0.1312"
    "E10" = "-4.59%"
    "F10" = "27-1-2023"
    "G10" = "0"
    "B11" = "WazirX"
    "C11" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D11" = "0.1896"
    "E11" = "-1.47%"
    "F11" = "27-1-2023"
    "G11" = "0"
    "B12" = "MandalaExchangeToken"
    "C12" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D12" = "0.09204"
    "E12" = "0.01%"
    "F12" = "27-1-2023"
    "G12" = "0"
    "B13" = "BitrueCoin"
    "C13" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D13" = "0.03518"
    "E13" = "0.92%"
    "F13" = "27-1-2023"
    "G13" = "0"
    "B14" = "BitMartToken"
    "C14" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D14" = "0.09894"
    "E14" = "0.42%"
    "F14" = "27-1-2023"
    "G14" = "0"
    "B15" = "BitForexToken"
    "C15" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D15" = "0.001430"
    "E15" = "-0.94%"
    "F15" = "27-1-2023"
    "G15" = "0"
    "B16" = "TigerCash"
    "C16" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D16" = "0.006283"
    "E16" = "7.02%"
    "F16" = "27-1-2023"
    "G16" = "0"
    "B17" = "LEO"
    "C17" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D17" = "3.658"
    "E17" = "3.53%"
    "F17" = "27-1-2023"
    "G17" = "0"
    "D18" = "3.169"
    "E18" = "4.73%"
    "F18" = "27-1-2023"
    "G18" = "0"
    "E19" = "0.34%"
    "F19" = "27-1-2023"
    "G19" = "0"
    "B20" = "MCDex"
    "C20" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "D20" = "5.207"
    "E20" = "4.58%"
    "F20" = "27-1-2023"
    "G20" = "0"
    "B21" = "ProBitToken"
    "C21" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "D21" = "0.1314"
    "E21" = "-1.57%"
    "F21" = "27-1-2023"
    "G21" = "0"
    "D22" = "0.2541"
    "E22" = "5.86%"
    "F22" = "27-1-2023"
    "G22" = "0"
    "D23" = "0.04448"
    "E23" = "-1.21%"
    "F23" = "27-1-2023"
    "G23" = "0"
    "D24" = "0.001240"
    "E24" = "2.67%"
    "F24" = "27-1-2023"
    "G24" = "0"
    "D25" = "0.004678"
    "E25" = "-5.33%"
    "F25" = "27-1-2023"
    "G25" = "0"
    "D26" = "0.0001306"
    "E26" = "5.59%"
    "F26" = "27-1-2023"
    "G26" = "0"
    "D27" = "0.0004460"
    "E27" = "0.67%"
    "F27" = "27-1-2023"
    "G27" = "0"
    "F28" = "27-1-2023"
    "G28" = "0"
    "F29" = "27-1-2023"
    "G29" = "0"
    "F30" = "27-1-2023"
    "G30" = "0"
    "F31" = "27-1-2023"
    "G31" = "0"
    "F32" = "27-1-2023"
    "G32" = "0"
    "F33" = "27-1-2023"
    "G33" = "0"
    "F34" = "27-1-2023"
    "G34" = "0"
    "F35" = "27-1-2023"
    "G35" = "0"
    "F36" = "27-1-2023"
    "G36" = "0"
    "F37" = "27-1-2023"
    "G37" = "0"
    "F38" = "27-1-2023"
    "G38" = "0"
    "D39" = "0.01925"
    "E39" = "-3.63%"
    "F39" = "27-1-2023"
    "G39" = "0"
    "D40" = "0.05114"
    "E40" = "3.47%"
    "F40" = "27-1-2023"
    "G40" = "0"
    "D41" = "0.007570"
    "E41" = "-0.32%"
    "F41" = "27-1-2023"
    "G41" = "0"
    "D42" = "0.01023"
    "E42" = "-7.39%"
    "F42" = "27-1-2023"
    "G42" = "0"
    "D43" = "0.1365"
    "E43" = "-1.69%"
    "F43" = "27-1-2023"
    "G43" = "0"
    "D44" = "0.002160"
    "E44" = "2.63%"
    "F44" = "27-1-2023"
    "G44" = "0"
    "D45" = "0.009889"
    "E45" = "-5.54%"
    "F45" = "27-1-2023"
    "G45" = "0"
    "D46" = "0.00006310"
    "E46" = "-3.34%"
    "F46" = "27-1-2023"
    "G46" = "0"
    "D47" = "0.00000000753"
    "E47" = "0.65%"
    "F47" = "27-1-2023"
    "G47" = "0"
    "D48" = "65.22"
    "E48" = "0.85%"
    "F48" = "27-1-2023"
    "G48" = "0"
    "D49" = "0.001666"
    "E49" = "40.29%"
    "F49" = "27-1-2023"
    "G49" = "0"
    "D50" = "0.00002109"
    "E50" = "0.65%"
    "F50" = "27-1-2023"
    "G50" = "0"
    "D51" = "0.0002008"
    "E51" = "0.65%"
    "F51" = "27-1-2023"
    "G51" = "0"
}

foreach ($ref in $updates.Keys) {
    Set-TextValue $ws $ref $updates[$ref]
}

Write-Host "Applied $($updates.Count) cell updates."
